# Remove the IPO record for 아이엠비디엑스 (listed 2024-04-03) from all three
# data sheets. Deleting the whole row (rather than just clearing cells) lets
# Excel shift the following rows up and drop the now-unused shared strings,
# matching how the workbook's row/string bookkeeping is expected to look
# afterwards.

$wb = $excel.ActiveWorkbook

# 01_리그테이블 — row 21 is the 아이엠비디엑스 entry
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows.Item(21).Delete()

# 02_통합집계_Rawdata — row 15 is the matching entry
$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows.Item(15).Delete()

# 03_IPO현황_Summary — row 9 is the matching entry
$ws3 = $wb.Worksheets.Item(3)
$ws3.Rows.Item(9).Delete()
